$d = $word.ActiveDocument

# --- 1. Remove the "Turns" row from the Blocks & Strikes table ---
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Title -eq "Blocks & Strikes") {
        $t.Rows.Item($t.Rows.Count).Delete()
    }
}

# --- 2. Remove "Front Stretch Kick", "Inside Crescent Kick", "Outside Crescent Kick" rows
#        from the Kicks table ---
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Title -eq "Kicks") {
        for ($r = $t.Rows.Count; $r -ge 1; $r--) {
            $label = $t.Rows.Item($r).Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7)
            if ($label -eq "Front Stretch Kick" -or $label -eq "Inside Crescent Kick" -or $label -eq "Outside Crescent Kick") {
                $t.Rows.Item($r).Delete()
            }
        }
    }
}

# --- 3. Remove the manual line break between "One Step Sparring" and the
#        "(Attacker: ...)" text, and tighten the table caption to match. ---
$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$searchRange.Find.Execute("One Step Sparring", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($searchRange.Find.Found) {
    $breakStart = $searchRange.End
    $breakRange = $d.Range($breakStart, $breakStart + 1)
    if ([int][char]($breakRange.Text[0]) -eq 11) {
        $breakRange.Delete()
    }
}

# Update the table caption (tblCaption) on the One Step Sparring table to match
# the now-contiguous heading text.
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Title -like "One Step Sparring *(Attacker*") {
        $t.Title = "One Step Sparring(Attacker: low block in front stance; Defender: Joon bee & nod, kick from back leg)"
    }
}

# --- 4. Append a "Revision: 05/01/24" paragraph at the end of the document ---
# (Use Content.Paragraphs rather than the Document-level Paragraphs collection:
#  after the row deletions above, Document.Paragraphs.Last can resolve stale.)
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$newPara = $d.Content.Paragraphs.Last
$newPara.Range.Text = "Revision: 05/01/24"
$newPara.Range.Style = "BodyText"
